# Insert a new row above row 2, shifting the existing data (rows 2-15) down
# to rows 3-16. Fill the newly inserted row with zeros, matching the
# "Add files via upload" commit which adds an extra "0 kgf/cm2" reference
# row to the top of the f_table lookup table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything below the header down by one row.
$ws.Rows("2:2").Insert()

# The freshly inserted row inherits the formatting of the row above it
# (the bordered/shaded "s=1" style); the new zero row should be plain,
# unstyled cells, so drop the inherited formatting before filling it in.
$ws.Range("A2:J2").ClearFormats()
$ws.Range("A2:J2").Value2 = 0

# Match the author's recorded selection at the time of saving.
$ws.Range("L4").Select() | Out-Null
